$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase the height of data rows 2 through 6
$ws.Range("A2:A6").EntireRow.RowHeight = 97.5
